# Apply updated crypto price/volume values to match the latest GitHub Actions scrape.
# Each target cell is briefly switched to a Text number format before the write so that
# numeric-looking strings (e.g. "17.00", "0.0422") are not silently coerced into floating
# point numbers (which would lose trailing zeros / precision). The style is then reset to
# "Normal" so the cell keeps the workbook default styling (no stray format is left behind).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "51.241.21"
Set-TextValue "E2" "  -0.82%  "

Set-TextValue "D3" "2.910.49"
Set-TextValue "E3" "  +0.67%  "

Set-TextValue "D5" "363.31"
Set-TextValue "E5" "  +3.04%  "

Set-TextValue "D6" "103.57"
Set-TextValue "E6" "  -4.65%  "

Set-TextValue "D7" "0.543"
Set-TextValue "E7" "  -2.94%  "

Set-TextValue "E8" "  +0.11%  "

Set-TextValue "D9" "0.589"
Set-TextValue "E9" "  -4.95%  "

Set-TextValue "D10" "36.80"
Set-TextValue "E10" "  -3.90%  "

Set-TextValue "D11" "0.139"
Set-TextValue "E11" "  +1.90%  "

Set-TextValue "E12" "  -3.14%  "

Set-TextValue "D13" "18.56"
Set-TextValue "E13" "  -3.82%  "

Set-TextValue "D14" "3.371.13"
Set-TextValue "E14" "  +0.57%  "

Set-TextValue "D15" "7.35"
Set-TextValue "E15" "  -3.69%  "

Set-TextValue "D16" "2.918.63"
Set-TextValue "E16" "  +1.71%  "

Set-TextValue "E17" "  -1.15%  "

Set-TextValue "D18" "51.203.00"
Set-TextValue "E18" "  -0.76%  "

Set-TextValue "E19" "  -1.53%  "

Set-TextValue "D20" "7.21"
Set-TextValue "E20" "  -3.27%  "

Set-TextValue "D21" "13.01"
Set-TextValue "E21" "  -4.75%  "

Set-TextValue "E22" "  -2.25%  "

Set-TextValue "D23" "68.30"
Set-TextValue "E23" "  -2.32%  "

Set-TextValue "D24" "259.54"
Set-TextValue "E24" "  -2.42%  "

Set-TextValue "E25" "  -2.42%  "

Set-TextValue "E26" "  -4.67%  "

Set-TextValue "E27" "  +0.06%  "

Set-TextValue "D28" "25.98"
Set-TextValue "E28" "  -2.16%  "

Set-TextValue "D29" "7.23"
Set-TextValue "E29" "  -3.03%  "

Set-TextValue "E30" "  +3.40%  "

Set-TextValue "D31" "6.16"
Set-TextValue "E31" "  +1.28%  "

Set-TextValue "D32" "9.94"
Set-TextValue "E32" "  -4.18%  "

Set-TextValue "D33" "2.14"
Set-TextValue "E33" "  -2.79%  "

Set-TextValue "D34" "34.89"
Set-TextValue "E34" "  -5.65%  "

Set-TextValue "D35" "50.56"
Set-TextValue "E35" "  -2.61%  "

Set-TextValue "E36" "  +0.14%  "

Set-TextValue "D37" "0.0422"
Set-TextValue "E37" "  -2.92%  "

Set-TextValue "E38" "  +4.65%  "

Set-TextValue "E39" "  +0.07%  "

Set-TextValue "D40" "17.00"
Set-TextValue "E40" "  -5.52%  "

Set-TextValue "E41" "  -5.61%  "

Set-TextValue "E42" "  -3.85%  "

Set-TextValue "D43" "22.36"
Set-TextValue "E43" "  -0.87%  "

Set-TextValue "D44" "119.23"
Set-TextValue "E44" "  +0.38%  "

Set-TextValue "D45" "2.14"
Set-TextValue "E45" "  -1.74%  "

Set-TextValue "D46" "2.072.17"
Set-TextValue "E46" "  -1.91%  "

Set-TextValue "E47" "  -5.95%  "

Set-TextValue "E48" "  -8.61%  "

Set-TextValue "D49" "3.202.66"

Set-TextValue "D50" "0.236"
Set-TextValue "E50" "  -4.37%  "

Set-TextValue "D51" "0.0308"
Set-TextValue "E51" "  -7.60%  "
